$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Window position (best effort - minor cosmetic, may not be exposed by host)
# ---------------------------------------------------------------------------
try { $excel.Windows.Item(1).Left = 1860 } catch {}

# ---------------------------------------------------------------------------
# PASS 1 - text labels, written in the exact order the original author must
# have typed them so the shared-string table comes out in the same order as
# the source file (the table is append-only / first-use-ordered).
# ---------------------------------------------------------------------------
$ws.Range("E14").Value = "shipping cost"

$ws.Range("B19").Value = "Description"
$ws.Range("C19").Value = "Hours"
$ws.Range("E19").Value = "amount"
$ws.Range("D19").Value = "Rate(`$/hr)"
$ws.Range("B19:E19").HorizontalAlignment = -4108

$ws.Range("B26").Value = "pcb layout work"
$ws.Range("B25").Value = "pcb Layout work"

$ws.Range("A31").Value = "Other expenses"
$ws.Range("A31").Font.Bold = $true

$ws.Range("A18").Value = "Hourly work"
$ws.Range("A18").Font.Bold = $true

$ws.Range("B32").Value = "single day shipping"
$ws.Range("C32").Value = "NA"

$ws.Range("A34").Value = "Rebtes"
$ws.Range("A34").Font.Bold = $true

$ws.Range("B35").Value = "Friends and family discount"
$ws.Range("B36").Value = "2 for 20 discount"

$ws.Range("F30").Value = "Hourly total"
$ws.Range("F33").Value = "other expenses total"
$ws.Range("F37").Value = "rebates total"
$ws.Range("F41").Value = "TOTAL"

# ---------------------------------------------------------------------------
# PASS 2 - remaining text cells that reuse already-known shared strings
# (order no longer matters for these).
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = "assembly"
$ws.Range("E12").Value = "assembly"

$ws.Range("B21").Value = "Schematic work "
$ws.Range("B22").Value = "Schematic work "
$ws.Range("B23").Value = "Schematic work "
$ws.Range("B24").Value = "Schematic work "
$ws.Range("B27").Value = "assembly"
$ws.Range("B28").Value = "assembly"
$ws.Range("B29").Value = "assembly"

$ws.Range("D32").Value = "NA"
$ws.Range("C35").Value = "NA"
$ws.Range("D35").Value = "NA"
$ws.Range("C36").Value = "NA"
$ws.Range("D36").Value = "NA"

# ---------------------------------------------------------------------------
# PASS 3 - numeric values, formulas and number formats (string table order
# is unaffected by these). NumberFormat is always applied to a row's C cell
# AFTER every other cell in that row has already been written - writing it
# earlier makes the host engine carry the same style onto the next *new*
# cell it sees, even several rows later.
# ---------------------------------------------------------------------------
$ws.Range("D14").Value = 66.64

$ws.Range("C21").Value = 0.11851851851851852
$ws.Range("D21").Value = 15
$ws.Range("E21").Formula = "=((C21-INT(C21))*24)*D21"
$ws.Range("C21").NumberFormat = $ws.Range("C4").NumberFormat

$ws.Range("C22").Value = 0.020833333333333332
$ws.Range("D22").Value = 15
$ws.Range("E22").Formula = "=((C22-INT(C22))*24)*D22"
$ws.Range("C22").NumberFormat = $ws.Range("C4").NumberFormat

$ws.Range("C23").Value = 0.022222222222222223
$ws.Range("D23").Value = 15
$ws.Range("E23").Formula = "=((C23-INT(C23))*24)*D23"
$ws.Range("C23").NumberFormat = $ws.Range("C4").NumberFormat

$ws.Range("C24").Value = 0.0625
$ws.Range("D24").Value = 15
$ws.Range("E24").Formula = "=((C24-INT(C24))*24)*D24"
$ws.Range("C24").NumberFormat = $ws.Range("C4").NumberFormat

$ws.Range("C25").Value = 0.3125
$ws.Range("D25").Value = 15
$ws.Range("E25").Formula = "=((C25-INT(C25))*24)*D25"
$ws.Range("C25").NumberFormat = $ws.Range("C4").NumberFormat

$ws.Range("C26").Value = 0.33333333333333331
$ws.Range("D26").Value = 15
$ws.Range("E26").Formula = "=((C26-INT(C26))*24)*D26"
$ws.Range("C26").NumberFormat = $ws.Range("C4").NumberFormat

$ws.Range("C27").Value = 0.09027777777777778
$ws.Range("D27").Value = 15
$ws.Range("E27").Formula = "=((C27-INT(C27))*24)*D27"
$ws.Range("C27").NumberFormat = $ws.Range("C10").NumberFormat

$ws.Range("C28").Value = 0.25
$ws.Range("D28").Value = 15
$ws.Range("E28").Formula = "=((C28-INT(C28))*24)*D28"
$ws.Range("C28").NumberFormat = $ws.Range("C10").NumberFormat

$ws.Range("C29").Value = 0.33333333333333331
$ws.Range("D29").Value = 15
$ws.Range("E29").Formula = "=((C29-INT(C29))*24)*D29"
$ws.Range("C29").NumberFormat = $ws.Range("C10").NumberFormat

$ws.Range("G30").Formula = "=SUM(E21:E29)"

$ws.Range("E32").Value = 66.64
$ws.Range("G33").Formula = "=SUM(E32)"

$ws.Range("E35").Value = -150
$ws.Range("E36").Value = -172.306
$ws.Range("G37").Formula = "=SUM(E35:E36)"

$ws.Range("G41").Formula = "=SUM(G30,G33,G37)"

$ws.Range("C11").Value = 0.25
$ws.Range("C12").Value = 0.33333333333333331
$ws.Range("C11").NumberFormat = $ws.Range("C10").NumberFormat
$ws.Range("C12").NumberFormat = $ws.Range("C10").NumberFormat

# ---------------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 15.92
$ws.Columns.Item(2).ColumnWidth = 22.76
$ws.Columns.Item(4).ColumnWidth = 10.42
$ws.Columns.Item(5).ColumnWidth = 12.25
$ws.Columns.Item(6).ColumnWidth = 16.76
$ws.Columns.Item(7).ColumnWidth = 11.76

# ---------------------------------------------------------------------------
# Page setup / view
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1

try { $wb.Windows.Item(1).ScrollRow = 16 } catch {}
$ws.Range("B29").Select()
